# The 2023 New York roster is reordered: the record for Jalen Brunson
# (row 5) and the record for RJ Barrett (row 6) trade places, while the
# sequential index in column A (3 / 4) stays put. Columns B..K hold the
# per-player data (No., Player, Pos, Ht, Wt, Birth Date, country, Exp,
# College, bbref url), so we swap those columns between the two rows.
# Row 14's College cell ("Villanova") is re-written too; its value does
# not change visibly, but this mirrors the shared-string reshuffle the
# row swap triggers upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

# Snapshot the current contents of row 5 and row 6 (use Value2, which
# reads back reliably for both text and numeric cells in this runtime).
$row5 = @{}
$row6 = @{}
foreach ($col in $cols) {
    $row5[$col] = $ws.Range($col + "5").Value2
    $row6[$col] = $ws.Range($col + "6").Value2
}

# Write each row's snapshot into the other row.
foreach ($col in $cols) {
    $ws.Range($col + "5").Value = $row6[$col]
    $ws.Range($col + "6").Value = $row5[$col]
}

# Row 14 (Josh Hart)'s College column keeps displaying "Villanova".
$ws.Range("J14").Value = "Villanova"
